$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 69.89967366666667
$ws.Range("H2").Value = 209.699021
$ws.Range("I2").Value = 0.6608367681537789
$ws.Range("J2").Value = 0.660836768153779
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 658.7242027828554
$ws.Range("R2").Value = 5928.517825045697
$ws.Range("S2").Value = 0.04397430088232663
$ws.Range("T2").Value = 0.04397430088232664
$ws.Range("G3").Value = 69.89967366666667
$ws.Range("H3").Value = 209.699021
$ws.Range("I3").Value = 0.6608367681537789
$ws.Range("J3").Value = 0.660836768153779
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 3536.403760159731
$ws.Range("R3").Value = 31827.63384143758
$ws.Range("S3").Value = 0.2360788966515605
$ws.Range("T3").Value = 0.2360788966515605
$ws.Range("G4").Value = 69.89967366666667
$ws.Range("H4").Value = 209.699021
$ws.Range("I4").Value = 0.6608367681537789
$ws.Range("J4").Value = 0.660836768153779
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 1876.16055136445
$ws.Range("R4").Value = 16885.44496228005
$ws.Range("S4").Value = 0.1252464206426748
$ws.Range("T4").Value = 0.1252464206426749
$ws.Range("G5").Value = 69.89967366666667
$ws.Range("H5").Value = 209.699021
$ws.Range("I5").Value = 0.6608367681537789
$ws.Range("J5").Value = 0.660836768153779
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 3827.883605258106
$ws.Range("R5").Value = 34450.95244732295
$ws.Range("S5").Value = 0.2555371499772169
$ws.Range("T5").Value = 0.255537149977217
$ws.Range("I6").Value = 0.1661491941864736
$ws.Range("J6").Value = 0.1661491941864736
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 165.6180478414756
$ws.Range("R6").Value = 1490.56243057328
$ws.Range("S6").Value = 0.01105612612464672
$ws.Range("T6").Value = 0.01105612612464673
$ws.Range("I7").Value = 0.1661491941864736
$ws.Range("J7").Value = 0.1661491941864736
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("S7").Value = 0.05935553276291208
$ws.Range("T7").Value = 0.05935553276291208
$ws.Range("I8").Value = 0.1661491941864736
$ws.Range("J8").Value = 0.1661491941864736
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 471.7088679017244
$ws.Range("R8").Value = 4245.37981111552
$ws.Range("S8").Value = 0.0314897609626922
$ws.Range("T8").Value = 0.0314897609626922
$ws.Range("I9").Value = 0.1661491941864736
$ws.Range("J9").Value = 0.1661491941864736
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 962.4158447328533
$ws.Range("R9").Value = 8661.742602595679
$ws.Range("S9").Value = 0.06424777433622261
$ws.Range("T9").Value = 0.06424777433622263
$ws.Range("G10").Value = 4.152730666666667
$ws.Range("H10").Value = 12.458192
$ws.Range("I10").Value = 0.0392602278210887
$ws.Range("J10").Value = 0.03926022782108871
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 39.13472058277156
$ws.Range("R10").Value = 352.212485244944
$ws.Range("S10").Value = 0.002612507587518945
$ws.Range("T10").Value = 0.002612507587518945
$ws.Range("G11").Value = 4.152730666666667
$ws.Range("H11").Value = 12.458192
$ws.Range("I11").Value = 0.0392602278210887
$ws.Range("J11").Value = 0.03926022782108871
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 210.0972947965831
$ws.Range("R11").Value = 1890.875653169248
$ws.Range("S11").Value = 0.0140254170363213
$ws.Range("T11").Value = 0.0140254170363213
$ws.Range("G12").Value = 4.152730666666667
$ws.Range("H12").Value = 12.458192
$ws.Range("I12").Value = 0.0392602278210887
$ws.Range("J12").Value = 0.03926022782108871
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 111.4624582425885
$ws.Range("R12").Value = 1003.162124183296
$ws.Range("S12").Value = 0.007440873821147722
$ws.Range("T12").Value = 0.007440873821147725
$ws.Range("G13").Value = 4.152730666666667
$ws.Range("H13").Value = 12.458192
$ws.Range("I13").Value = 0.0392602278210887
$ws.Range("J13").Value = 0.03926022782108871
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 227.4140750898293
$ws.Range("R13").Value = 2046.726675808464
$ws.Range("S13").Value = 0.01518142937610073
$ws.Range("T13").Value = 0.01518142937610074
$ws.Range("G14").Value = 14.14774133333333
$ws.Range("H14").Value = 42.443224
$ws.Range("I14").Value = 0.1337538098386587
$ws.Range("J14").Value = 0.1337538098386588
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 133.3262251755298
$ws.Range("R14").Value = 1199.936026579768
$ws.Range("S14").Value = 0.008900428307636146
$ws.Range("T14").Value = 0.008900428307636148
$ws.Range("G15").Value = 14.14774133333333
$ws.Range("H15").Value = 42.443224
$ws.Range("I15").Value = 0.1337538098386587
$ws.Range("J15").Value = 0.1337538098386588
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 715.7705182939396
$ws.Range("R15").Value = 6441.934664645456
$ws.Range("S15").Value = 0.04778252871411846
$ws.Range("T15").Value = 0.04778252871411847
$ws.Range("G16").Value = 14.14774133333333
$ws.Range("H16").Value = 42.443224
$ws.Range("I16").Value = 0.1337538098386587
$ws.Range("J16").Value = 0.1337538098386588
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 379.7361673973903
$ws.Range("R16").Value = 3417.625506576512
$ws.Range("S16").Value = 0.0253499604394208
$ws.Range("T16").Value = 0.02534996043942081
$ws.Range("G17").Value = 14.14774133333333
$ws.Range("H17").Value = 42.443224
$ws.Range("I17").Value = 0.1337538098386587
$ws.Range("J17").Value = 0.1337538098386588
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 774.7662365285787
$ws.Range("R17").Value = 6972.896128757207
$ws.Range("S17").Value = 0.05172089237748331
$ws.Range("T17").Value = 0.05172089237748333
